# The commit replaces the TestData.xlsx with an updated version: two of the
# random-looking id values on Sheet1 (row 2) were regenerated, and the sheet's
# active selection was left on D10 when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Post id" (B2) and "User id" (C2) values. These cells are
# formatted as Text, and the values keep their leading backtick, so use
# single-quoted strings to avoid PowerShell interpreting the backtick as an
# escape character.
$ws.Range("B2").Value = '`233215'
$ws.Range("C2").Value = '`8078407'

# Leave the selection on D10, as in the saved workbook.
$ws.Range("D10").Select() | Out-Null
